$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header counts in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: remove stray D2 value, add new B2/C2 values
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 42.117013549239076
$ws.Range("C2").Value = 21.142373404282935

# Row 3: remove stray B3 value, update C3 value
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 41.332072333428492

# Update the selected range to match the new data extent
$ws.Range("B1:E3").Select()
